# "Colocando header nos gráficos"
# Adds a header label to column A (row 1) on each data sheet, fixes
# accented Portuguese text in the row labels (which also drops their
# bold/border header style, matching the author's edit), removes the
# obsolete "Teto" row from the Emissoes sheet, and updates the Custo
# Total sheet (new header row label + refreshed numbers).

$wb = $excel.ActiveWorkbook

function Set-HeaderCell {
    param($ws, [string]$cell, [string]$text, [string]$styleSourceCell)
    # Clone formatting (font/border/alignment => same style index) from an
    # existing header cell, then set the text.
    $ws.Range($styleSourceCell).Copy() | Out-Null
    $ws.Range($cell).PasteSpecial(-4122) | Out-Null
    $ws.Range($cell).Value = $text
}

# ---------------------------------------------------------------------
# Sheets 1-4: "Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio
# (MWMed)", "Atendimento a Ponta(MW)", "Potencia Incremental - SIN(MW)"
# All four share the same row layout (rows 2-12 = technologies).
# ---------------------------------------------------------------------
$techSheets = @(1, 2, 3, 4)
foreach ($idx in $techSheets) {
    $ws = $wb.Worksheets.Item($idx)

    Set-HeaderCell $ws "A1" "Fonte/Tecnologia" "B1"

    $labels = @{
        2  = "Hidro"
        3  = "Gás Natural"
        4  = "Carvão"
        5  = "Nuclear"
        6  = "Óleos Comb"
        7  = "Biomassa"
        8  = "Eólica"
        9  = "Solar"
        10 = "Outros"
        11 = "Pot. Compl."
        12 = "GD"
    }

    foreach ($row in 2..12) {
        $cell = "A$row"
        $ws.Range($cell).ClearFormats()
        $ws.Range($cell).Value = $labels[$row]
    }
}

# ---------------------------------------------------------------------
# Sheet 5: "Emissoes Totais (MtCO2eq)"
# Add header, rename/unstyle rows 2-3, remove the row 4 "Teto".
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

Set-HeaderCell $ws5 "A1" "Período" "B1"

$ws5.Range("A2").ClearFormats()
$ws5.Range("A2").Value = "P.Médio"

$ws5.Range("A3").ClearFormats()
$ws5.Range("A3").Value = "P.Crítico"

$ws5.Rows.Item(4).Delete()

# ---------------------------------------------------------------------
# Sheet 6: "Custo Total (bilhões de R$)"
# Add header, change B1 label, rename/unstyle rows 2-3 and refresh values.
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

Set-HeaderCell $ws6 "A1" "Tipo Expansão" "B1"

$ws6.Range("B1").Value = "2015"

$ws6.Range("A2").ClearFormats()
$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 577

$ws6.Range("A3").ClearFormats()
$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99
